# Update the "required share peak" (B3) parameter on the CDRDfRCP sheet
# from 0.05 to 0, then return focus to the About sheet (the workbook's
# active tab), leaving the cursor on CDRDfRCP parked at B4 (the cell
# just below the edited one, as if Enter was pressed after the edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDRDfRCP")

$ws.Range("B3").Value = 0
$ws.Range("B4").Select()

$wb.Worksheets.Item("About").Activate()
